$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 37.88856266666667
$ws.Range("H2").Value = 113.665688
$ws.Range("I2").Value = 0.9468476050819132
$ws.Range("J2").Value = 0.9540268599258594
$ws.Range("M2").Value = 1.038069333333333
$ws.Range("N2").Value = 3.114208
$ws.Range("O2").Value = 0.01303513967359888
$ws.Range("P2").Value = 0.01321546873134007
$ws.Range("Q2").Value = 39.33095498834489
$ws.Range("R2").Value = 353.978594895104
$ws.Range("S2").Value = 0.01234229078185533
$ws.Range("T2").Value = 0.01260791213620875

$ws.Range("G3").Value = 37.88856266666667
$ws.Range("H3").Value = 113.665688
$ws.Range("I3").Value = 0.9468476050819132
$ws.Range("J3").Value = 0.9540268599258594
$ws.Range("M3").Value = 75.23900966666666
$ws.Range("O3").Value = 0.9447837137804441
$ws.Range("P3").Value = 0.9578539194814477
$ws.Range("Q3").Value = 2850.697932733439
$ws.Range("R3").Value = 25656.28139460095
$ws.Range("S3").Value = 0.8945661967134093
$ws.Range("T3").Value = 0.9138183670705625

$ws.Range("G4").Value = 37.88856266666667
$ws.Range("H4").Value = 113.665688
$ws.Range("I4").Value = 0.9468476050819132
$ws.Range("J4").Value = 0.9540268599258594
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.062446
$ws.Range("N4").Value = 0.187338
$ws.Range("O4").Value = 0.0007841406213626924
$ws.Range("P4").Value = 0.0007949884789942698
$ws.Range("Q4").Value = 2.365989184282667
$ws.Range("R4").Value = 21.293902658544
$ws.Range("S4").Value = 0.0007424616693847086
$ws.Range("T4").Value = 0.0007584403622921382

$ws.Range("G5").Value = 37.88856266666667
$ws.Range("H5").Value = 113.665688
$ws.Range("I5").Value = 0.9468476050819132
$ws.Range("J5").Value = 0.9540268599258594
$ws.Range("M5").Value = 3.259981
$ws.Range("N5").Value = 6.519962
$ws.Range("O5").Value = 0.04093590505349536
$ws.Range("P5").Value = 0.02766814353457621
$ws.Range("Q5").Value = 123.5159944106426
$ws.Range("R5").Value = 741.0959664638559
$ws.Range("S5").Value = 0.03876006366176267
$ws.Range("T5").Value = 0.02639615209626971

$ws.Range("G6").Value = 37.88856266666667
$ws.Range("H6").Value = 113.665688
$ws.Range("I6").Value = 0.9468476050819132
$ws.Range("J6").Value = 0.9540268599258594
$ws.Range("M6").Value = 0.03672033333333333
$ws.Range("N6").Value = 0.110161
$ws.Range("O6").Value = 0.0004611008710989525
$ws.Range("P6").Value = 0.0004674797736416945
$ws.Range("Q6").Value = 1.391280650640889
$ws.Range("R6").Value = 12.521525855768
$ws.Range("S6").Value = 0.0004365922555012271
$ws.Range("T6").Value = 0.0004459882605262373

$ws.Range("I7").Value = 0.02931771140176381
$ws.Range("J7").Value = 0.02954000622562442
$ws.Range("M7").Value = 1.038069333333333
$ws.Range("N7").Value = 3.114208
$ws.Range("O7").Value = 0.01303513967359888
$ws.Range("P7").Value = 0.01321546873134007
$ws.Range("Q7").Value = 1.217823841255111
$ws.Range("R7").Value = 10.960414571296
$ws.Range("S7").Value = 0.0003821604630322537
$ws.Range("T7").Value = 0.0003903850285983306

$ws.Range("I8").Value = 0.02931771140176381
$ws.Range("J8").Value = 0.02954000622562442
$ws.Range("M8").Value = 75.23900966666666
$ws.Range("O8").Value = 0.9447837137804441
$ws.Range("P8").Value = 0.9578539194814477
$ws.Range("Q8").Value = 88.26757213823588
$ws.Range("R8").Value = 794.4081492441228
$ws.Range("S8").Value = 0.02769889625770168
$ws.Range("T8").Value = 0.02829501074472072

$ws.Range("I9").Value = 0.02931771140176381
$ws.Range("J9").Value = 0.02954000622562442
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.062446
$ws.Range("N9").Value = 0.187338
$ws.Range("O9").Value = 0.0007841406213626924
$ws.Range("P9").Value = 0.0007949884789942698
$ws.Range("Q9").Value = 0.07325929506733332
$ws.Range("R9").Value = 0.659333655606
$ws.Range("S9").Value = 0.00002298920843551116
$ws.Range("T9").Value = 0.00002348396461879041

$ws.Range("I10").Value = 0.02931771140176381
$ws.Range("J10").Value = 0.02954000622562442
$ws.Range("M10").Value = 3.259981
$ws.Range("N10").Value = 6.519962
$ws.Range("O10").Value = 0.04093590505349536
$ws.Range("P10").Value = 0.02766814353457621
$ws.Range("Q10").Value = 3.824486916582333
$ws.Range("R10").Value = 22.946921499494
$ws.Range("S10").Value = 0.001200147050328381
$ws.Range("T10").Value = 0.0008173171322628511

$ws.Range("I11").Value = 0.02931771140176381
$ws.Range("J11").Value = 0.02954000622562442
$ws.Range("M11").Value = 0.03672033333333333
$ws.Range("N11").Value = 0.110161
$ws.Range("O11").Value = 0.0004611008710989525
$ws.Range("P11").Value = 0.0004674797736416945
$ws.Range("Q11").Value = 0.04307891193411111
$ws.Range("R11").Value = 0.387710207407
$ws.Range("S11").Value = 0.00001351842226598098
$ws.Range("T11").Value = 0.00001380935542372915

$ws.Range("G12").Value = 0.05038133333333333
$ws.Range("H12").Value = 0.151144
$ws.Range("I12").Value = 0.00125904604054744
$ws.Range("J12").Value = 0.001268592468438093
$ws.Range("M12").Value = 1.038069333333333
$ws.Range("N12").Value = 3.114208
$ws.Range("O12").Value = 0.01303513967359888
$ws.Range("P12").Value = 0.01321546873134007
$ws.Range("Q12").Value = 0.05229931710577777
$ws.Range("R12").Value = 0.470693853952
$ws.Range("S12").Value = 0.00001641184099402753
$ws.Range("T12").Value = 0.00001676504409945713

$ws.Range("G13").Value = 0.05038133333333333
$ws.Range("H13").Value = 0.151144
$ws.Range("I13").Value = 0.00125904604054744
$ws.Range("J13").Value = 0.001268592468438093
$ws.Range("M13").Value = 75.23900966666666
$ws.Range("O13").Value = 0.9447837137804441
$ws.Range("P13").Value = 0.9578539194814477
$ws.Range("Q13").Value = 3.790641625686222
$ws.Range("R13").Value = 34.115774631176
$ws.Range("S13").Value = 0.001189526194008974
$ws.Range("T13").Value = 0.001215126268118072

$ws.Range("G14").Value = 0.05038133333333333
$ws.Range("H14").Value = 0.151144
$ws.Range("I14").Value = 0.00125904604054744
$ws.Range("J14").Value = 0.001268592468438093
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.062446
$ws.Range("N14").Value = 0.187338
$ws.Range("O14").Value = 0.0007841406213626924
$ws.Range("P14").Value = 0.0007949884789942698
$ws.Range("Q14").Value = 0.003146112741333333
$ws.Range("R14").Value = 0.028315014672
$ws.Range("S14").Value = 0.0000009872691445591075
$ws.Range("T14").Value = 0.000001008516396947186

$ws.Range("G15").Value = 0.05038133333333333
$ws.Range("H15").Value = 0.151144
$ws.Range("I15").Value = 0.00125904604054744
$ws.Range("J15").Value = 0.001268592468438093
$ws.Range("M15").Value = 3.259981
$ws.Range("N15").Value = 6.519962
$ws.Range("O15").Value = 0.04093590505349536
$ws.Range("P15").Value = 0.02766814353457621
$ws.Range("Q15").Value = 0.1642421894213333
$ws.Range("R15").Value = 0.9854531365279999
$ws.Range("S15").Value = 0.00005154018917382928
$ws.Range("T15").Value = 0.00003509959850362749

$ws.Range("G16").Value = 0.05038133333333333
$ws.Range("H16").Value = 0.151144
$ws.Range("I16").Value = 0.00125904604054744
$ws.Range("J16").Value = 0.001268592468438093
$ws.Range("M16").Value = 0.03672033333333333
$ws.Range("N16").Value = 0.110161
$ws.Range("O16").Value = 0.0004611008710989525
$ws.Range("P16").Value = 0.0004674797736416945
$ws.Range("Q16").Value = 0.001850019353777778
$ws.Range("R16").Value = 0.016650174184
$ws.Range("S16").Value = 0.0000005805472260501117
$ws.Range("T16").Value = 0.000000593041319988998

$ws.Range("G17").Value = 0.903375
$ws.Range("H17").Value = 1.80675
$ws.Range("I17").Value = 0.02257563747577563
$ws.Range("J17").Value = 0.0151645413800781
$ws.Range("M17").Value = 1.038069333333333
$ws.Range("N17").Value = 3.114208
$ws.Range("O17").Value = 0.01303513967359888
$ws.Range("P17").Value = 0.01321546873134007
$ws.Range("Q17").Value = 0.937765884
$ws.Range("R17").Value = 5.626595303999999
$ws.Range("S17").Value = 0.0002942765877172686
$ws.Range("T17").Value = 0.0002004065224335348

$ws.Range("G18").Value = 0.903375
$ws.Range("H18").Value = 1.80675
$ws.Range("I18").Value = 0.02257563747577563
$ws.Range("J18").Value = 0.0151645413800781
$ws.Range("M18").Value = 75.23900966666666
$ws.Range("O18").Value = 0.9447837137804441
$ws.Range("P18").Value = 0.9578539194814477
$ws.Range("Q18").Value = 67.96904035762499
$ws.Range("R18").Value = 407.81424214575
$ws.Range("S18").Value = 0.02132909461532427
$ws.Range("T18").Value = 0.01452541539804641

$ws.Range("G19").Value = 0.903375
$ws.Range("H19").Value = 1.80675
$ws.Range("I19").Value = 0.02257563747577563
$ws.Range("J19").Value = 0.0151645413800781
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.062446
$ws.Range("N19").Value = 0.187338
$ws.Range("O19").Value = 0.0007841406213626924
$ws.Range("P19").Value = 0.0007949884789942698
$ws.Range("Q19").Value = 0.05641215525000001
$ws.Range("R19").Value = 0.3384729315
$ws.Range("S19").Value = 0.00001770247439791359
$ws.Range("T19").Value = 0.00001205563568639395

$ws.Range("G20").Value = 0.903375
$ws.Range("H20").Value = 1.80675
$ws.Range("I20").Value = 0.02257563747577563
$ws.Range("J20").Value = 0.0151645413800781
$ws.Range("M20").Value = 3.259981
$ws.Range("N20").Value = 6.519962
$ws.Range("O20").Value = 0.04093590505349536
$ws.Range("P20").Value = 0.02766814353457621
$ws.Range("Q20").Value = 2.944985335875
$ws.Range("R20").Value = 11.7799413435
$ws.Range("S20").Value = 0.0009241541522304827
$ws.Range("T20").Value = 0.0004195747075400212

$ws.Range("G21").Value = 0.903375
$ws.Range("H21").Value = 1.80675
$ws.Range("I21").Value = 0.02257563747577563
$ws.Range("J21").Value = 0.0151645413800781
$ws.Range("M21").Value = 0.03672033333333333
$ws.Range("N21").Value = 0.110161
$ws.Range("O21").Value = 0.0004611008710989525
$ws.Range("P21").Value = 0.0004674797736416945
$ws.Range("Q21").Value = 0.033172231125
$ws.Range("R21").Value = 0.19903338675
$ws.Range("S21").Value = 0.0000104096461056943
$ws.Range("T21").Value = 0.000007089116371739019
